# Update diagnostics results 5th May
# Adds a new worksheet "2017-05-04" (after the existing "2017-05-02" sheet)
# containing the 5-May data-entry performance snapshot, mirroring the
# layout/header of the existing sheet.

$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after the existing "2017-05-02" sheet.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$ws.Name = "2017-05-04"

# Header row (reuses the same shared-string header labels as sheet 1).
$ws.Range("A1").Value = "CreatedBy"
$ws.Range("B1").Value = "DataCompleted"
$ws.Range("C1").Value = "DataIncomplete"
$ws.Range("D1").Value = "Total"

$data = @(
    @('Default Administrator', 0, 1, 1),
    @('Dr Masliyana', 1, 0, 1),
    @('Fadzli', 0, 1, 1),
    @('Syahirah', 0, 1, 1),
    @('Aisyah', 2, 0, 2),
    @('DR ANG', 2, 0, 2),
    @('Dr Saravanan', 2, 0, 2),
    @('Nabilah Iffah', 2, 0, 2),
    @('Sree Durga', 2, 0, 2),
    @('Aminiril anisah', 3, 0, 3),
    @('Adlan', 4, 0, 4),
    @('Aimi Nadiah Jamel', 4, 0, 4),
    @('Danial', 4, 0, 4),
    @('eliza', 4, 0, 4),
    @('Hui Yi', 4, 0, 4),
    @('Jia yi', 4, 0, 4),
    @('Yi Shin', 4, 0, 4),
    @('Yvonne', 4, 0, 4),
    @('Delina', 5, 0, 5),
    @('Ling Kuok Wei', 5, 0, 5),
    @('Dr saravanan', 7, 0, 7),
    @('Helmi', 8, 0, 8),
    @('Hooi Fan', 5, 3, 8),
    @('Jennifer Kaur', 8, 0, 8),
    @('Munirah', 8, 0, 8),
    @('Natrah', 7, 1, 8),
    @('Hui Che', 9, 0, 9),
    @('Hoong Ping', 10, 0, 10),
    @('Michelle', 9, 1, 10),
    @('Philip', 10, 0, 10),
    @('Thivashini', 9, 1, 10),
    @('Suhayl', 11, 0, 11),
    @('Annas', 12, 0, 12),
    @('Noor Amalina', 15, 0, 15),
    @('Yhyviyaa', 16, 0, 16),
    @('Musfirah', 16, 1, 17),
    @('Hadi', 18, 0, 18),
    @('Dr Richard', 19, 0, 19),
    @('Izzat', 19, 0, 19),
    @('Nurjannah', 19, 0, 19),
    @('Tan khai shin', 19, 0, 19),
    @('Eliza', 22, 0, 22),
    @('Geetha Krishnan', 21, 1, 22),
    @('Mardhiah', 21, 1, 22),
    @('Nursyuhaida', 22, 1, 23),
    @('Izzati', 24, 0, 24),
    @('Pui Yee', 23, 1, 24),
    @('Total', 443, 14, 457)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

# Keep the originally-active sheet ("2017-05-02") selected, since the
# workbook's bookViews/activeTab is unchanged by this edit.
$firstSheet.Activate()
